# Scheduled runner update: refresh Leve profit calculations (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the ALC, ARM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15 (G15=44146)
$ws.Range("H15").Value = 1542.8823
$ws.Range("I15").Value = 1542.8823
$ws.Range("K15").Value = 4628.6469
$ws.Range("M15").Value = -4459.6469

# Row 33 (G33=5512)
$ws.Range("H33").Value = 309.23077
$ws.Range("I33").Value = 302.9565
$ws.Range("K33").Value = 302.9565
$ws.Range("M33").Value = -73.95650000000001

# Row 107 (G107=27766)
$ws.Range("H107").Value = 993.8182
$ws.Range("I107").Value = 696.8
$ws.Range("J107").Value = 1630.2858
$ws.Range("K107").Value = 696.8
$ws.Range("L107").Value = 1630.2858
$ws.Range("M107").Value = 1223.2
$ws.Range("N107").Value = -5470.2858

# Row 130 (G130=34691)
$ws.Range("H130").Value = 92498.5
$ws.Range("J130").Value = 92498.5
$ws.Range("L130").Value = 92498.5
$ws.Range("N130").Value = -102538.5

# Row 135 (G135=44047)
$ws.Range("H135").Value = 1207.2424
$ws.Range("I135").Value = 1080.9565
$ws.Range("J135").Value = 1497.7
$ws.Range("K135").Value = 9728.6085
$ws.Range("L135").Value = 13479.3
$ws.Range("M135").Value = -7193.6085
$ws.Range("N135").Value = -18549.3

# Row 138 (G138=44169)
$ws.Range("H138").Value = 4344.533
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 4344.533
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 13033.599
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -23313.599


$ws = $wb.Worksheets.Item("ARM")
# Row 60 (G60=3883)
$ws.Range("H60").Value = 67410.55499999999
$ws.Range("I60").Value = 63337
$ws.Range("K60").Value = 63337
$ws.Range("M60").Value = -62604

# Row 61 (G61=43999)
$ws.Range("H61").Value = 336205.44
$ws.Range("I61").Value = 2213.4814
$ws.Range("J61").Value = 3342133
$ws.Range("K61").Value = 2213.4814
$ws.Range("L61").Value = 3342133
$ws.Range("M61").Value = -2001.4814
$ws.Range("N61").Value = -3342557

# Row 74 (G74=44000)
$ws.Range("H74").Value = 1175.6666
$ws.Range("I74").Value = 870.3182
$ws.Range("J74").Value = 2519.2
$ws.Range("K74").Value = 870.3182
$ws.Range("L74").Value = 2519.2
$ws.Range("M74").Value = 3.681799999999953
$ws.Range("N74").Value = -4267.2

# Row 77 (G77=44000)
$ws.Range("H77").Value = 1175.6666
$ws.Range("I77").Value = 870.3182
$ws.Range("J77").Value = 2519.2
$ws.Range("K77").Value = 4351.591
$ws.Range("L77").Value = 12596
$ws.Range("M77").Value = 16.40899999999965
$ws.Range("N77").Value = -21332

# Row 122 (G122=36168)
$ws.Range("H122").Value = 5704
$ws.Range("I122").Value = 4460.8125
$ws.Range("K122").Value = 13382.4375
$ws.Range("M122").Value = -10932.4375

# Row 132 (G132=43997)
$ws.Range("H132").Value = 4548887.5
$ws.Range("I132").Value = 2923.5715
$ws.Range("J132").Value = 12504324
$ws.Range("K132").Value = 8770.7145
$ws.Range("L132").Value = 37512972
$ws.Range("M132").Value = -6240.7145
$ws.Range("N132").Value = -37518032

# Row 136 (G136=43999)
$ws.Range("H136").Value = 336205.44
$ws.Range("I136").Value = 2213.4814
$ws.Range("J136").Value = 3342133
$ws.Range("K136").Value = 6640.4442
$ws.Range("L136").Value = 10026399
$ws.Range("M136").Value = -4090.4442
$ws.Range("N136").Value = -10031499

# Row 138 (G138=42350)
$ws.Range("H138").Value = 95429
$ws.Range("J138").Value = 95429
$ws.Range("L138").Value = 95429
$ws.Range("N138").Value = -105709


$ws = $wb.Worksheets.Item("CRP")
# Row 70 (G70=12011)
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# Row 73 (G73=12011)
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# Row 105 (G105=19928)
$ws.Range("H105").Value = 1353.8334
$ws.Range("I105").Value = 1176.5714
$ws.Range("J105").Value = 1602
$ws.Range("K105").Value = 1176.5714
$ws.Range("L105").Value = 1602
$ws.Range("M105").Value = 570.4286
$ws.Range("N105").Value = -5096

# Row 122 (G122=36196)
$ws.Range("H122").Value = 3014.72
$ws.Range("I122").Value = 3172.8333
$ws.Range("J122").Value = 2868.7693
$ws.Range("K122").Value = 9518.499899999999
$ws.Range("L122").Value = 8606.3079
$ws.Range("M122").Value = -7068.499899999999
$ws.Range("N122").Value = -13506.3079


$ws = $wb.Worksheets.Item("CUL")
# Row 12 (G12=4854)
$ws.Range("H12").Value = 3604.3635
$ws.Range("J12").Value = 4947.125
$ws.Range("L12").Value = 14841.375
$ws.Range("N12").Value = -15187.375

# Row 14 (G14=12886)
$ws.Range("H14").Value = 13479.538
$ws.Range("I14").Value = 13479.538
$ws.Range("K14").Value = 40438.614
$ws.Range("M14").Value = -40265.614

# Row 18 (G18=36056)
$ws.Range("H18").Value = 447.25
$ws.Range("I18").Value = 368.2857
$ws.Range("K18").Value = 1104.8571
$ws.Range("M18").Value = -935.8571000000002

# Row 25 (G25=4709)
$ws.Range("H25").Value = 11565.6
$ws.Range("I25").Value = 6000
$ws.Range("K25").Value = 18000
$ws.Range("M25").Value = -17831

# Row 30 (G30=4709)
$ws.Range("H30").Value = 11565.6
$ws.Range("I30").Value = 6000
$ws.Range("K30").Value = 18000
$ws.Range("M30").Value = -17898

# Row 70 (G70=12867)
$ws.Range("H70").Value = 15546.889
$ws.Range("I70").Value = 11265
$ws.Range("J70").Value = 24110.666
$ws.Range("K70").Value = 33795
$ws.Range("L70").Value = 72331.99800000001
$ws.Range("M70").Value = -33480
$ws.Range("N70").Value = -72961.99800000001

# Row 73 (G73=12867)
$ws.Range("H73").Value = 15546.889
$ws.Range("I73").Value = 11265
$ws.Range("J73").Value = 24110.666
$ws.Range("K73").Value = 33795
$ws.Range("L73").Value = 72331.99800000001
$ws.Range("M73").Value = -32703
$ws.Range("N73").Value = -74515.99800000001

# Row 82 (G82=12856)
$ws.Range("H82").Value = 12867
$ws.Range("I82").Value = 7001
$ws.Range("K82").Value = 21003
$ws.Range("M82").Value = -20597

# Row 85 (G85=12856)
$ws.Range("H85").Value = 12867
$ws.Range("I85").Value = 7001
$ws.Range("K85").Value = 21003
$ws.Range("M85").Value = -19599

# Row 129 (G129=36054)
$ws.Range("H129").Value = 1091715.6
$ws.Range("I129").Value = 2922
$ws.Range("J129").Value = 1791654.4
$ws.Range("K129").Value = 8766
$ws.Range("L129").Value = 5374963.199999999
$ws.Range("M129").Value = -3766
$ws.Range("N129").Value = -5384963.199999999

# Row 140 (G140=44097)
$ws.Range("H140").Value = 2877.3333
$ws.Range("I140").Value = 1199.7391
$ws.Range("J140").Value = 6735.8
$ws.Range("K140").Value = 3599.2173
$ws.Range("L140").Value = 20207.4
$ws.Range("M140").Value = 1580.7827
$ws.Range("N140").Value = -30567.4


$ws = $wb.Worksheets.Item("GSM")
# Row 64 (G64=10640)
$ws.Range("H64").Value = 60000.555
$ws.Range("J64").Value = 60000.555
$ws.Range("L64").Value = 60000.555
$ws.Range("N64").Value = -60496.555

# Row 67 (G67=10640)
$ws.Range("H67").Value = 60000.555
$ws.Range("J67").Value = 60000.555
$ws.Range("L67").Value = 60000.555
$ws.Range("N67").Value = -61716.555

# Row 70 (G70=14146)
$ws.Range("H70").Value = 8270.923000000001
$ws.Range("I70").Value = 7889.4
$ws.Range("J70").Value = 8509.375
$ws.Range("K70").Value = 7889.4
$ws.Range("L70").Value = 8509.375
$ws.Range("M70").Value = -7619.4
$ws.Range("N70").Value = -9049.375

# Row 73 (G73=14146)
$ws.Range("H73").Value = 8270.923000000001
$ws.Range("I73").Value = 7889.4
$ws.Range("J73").Value = 8509.375
$ws.Range("K73").Value = 7889.4
$ws.Range("L73").Value = 8509.375
$ws.Range("M73").Value = -6953.4
$ws.Range("N73").Value = -10381.375


$ws = $wb.Worksheets.Item("LTW")
# Row 107 (G107=38752)
$ws.Range("H107").Value = 16679.334
$ws.Range("I107").Value = 16679.334
$ws.Range("K107").Value = 16679.334
$ws.Range("M107").Value = -14759.334

# Row 132 (G132=44058)
$ws.Range("H132").Value = 2070.8438
$ws.Range("I132").Value = 2009.2333
$ws.Range("J132").Value = 2995
$ws.Range("K132").Value = 6027.699900000001
$ws.Range("L132").Value = 8985
$ws.Range("M132").Value = -3497.699900000001
$ws.Range("N132").Value = -14045

# Row 136 (G136=44060)
$ws.Range("H136").Value = 2221.3823
$ws.Range("I136").Value = 2152.8386
$ws.Range("J136").Value = 2929.6667
$ws.Range("K136").Value = 6458.5158
$ws.Range("L136").Value = 8789.000100000001
$ws.Range("M136").Value = -3908.5158
$ws.Range("N136").Value = -13889.0001


$ws = $wb.Worksheets.Item("WVR")
# Row 136 (G136=44031)
$ws.Range("H136").Value = 245063.9
$ws.Range("I136").Value = 7112.3057
$ws.Range("J136").Value = 1672773.5
$ws.Range("K136").Value = 21336.9171
$ws.Range("L136").Value = 5018320.5
$ws.Range("M136").Value = -18786.9171
$ws.Range("N136").Value = -5023420.5

# Row 138 (G138=42347)
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

Write-Host "Applied 209 cell updates across 8 sheets."